$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp (row 1)
$ws.Range("A1").Value = "Datos actualizados a 19 de Agosto de 2020 a las 08:19"

# India (row 6) - updated case counts
$ws.Range("B6").Value = 2768670
$ws.Range("C6").Value = 2044
$ws.Range("D6").Value = 2038585
$ws.Range("E6").Value = 677059
$ws.Range("G6").Value = 12
$ws.Range("H6").Value = 53026

# Alemania (row 22) - updated case counts
$ws.Range("D22").Value = 203900
$ws.Range("E22").Value = 14900

# Israel (row 33) - updated case counts
$ws.Range("B33").Value = 96753
$ws.Range("C33").Value = 344
$ws.Range("D33").Value = 72494
$ws.Range("E33").Value = 23548
$ws.Range("G33").Value = 3
$ws.Range("H33").Value = 711

# El Salvador (row 73) - updated case counts
$ws.Range("D73").Value = 11205
$ws.Range("E73").Value = 11632

# Rows 147/148: Georgia and Republica de Chipre swap order (Georgia's
# numbers updated, Republica de Chipre stays the same but moves down)
$ws.Range("A147").Value = "Georgia"
$ws.Range("B147").Value = 1361
$ws.Range("C147").Value = 10
$ws.Range("D147").Value = 1098
$ws.Range("E147").Value = 246
$ws.Range("F147").Value = 0
$ws.Range("G147").Value = 0
$ws.Range("H147").Value = 17

$ws.Range("A148").Value = "Republica de Chipre"
$ws.Range("B148").Value = 1359
$ws.Range("C148").Value = 0
$ws.Range("D148").Value = 878
$ws.Range("E148").Value = 461
$ws.Range("F148").Value = 0
$ws.Range("G148").Value = 0
$ws.Range("H148").Value = 20

# Rows 213/214: Montserrat and Islas Malvinas swap order (values
# unchanged, only the row order/country name swaps)
$ws.Range("A213").Value = "Montserrat"
$ws.Range("B213").Value = 13
$ws.Range("C213").Value = 0
$ws.Range("D213").Value = 12
$ws.Range("E213").Value = 0
$ws.Range("F213").Value = 0
$ws.Range("G213").Value = 0
$ws.Range("H213").Value = 1

$ws.Range("A214").Value = "Islas Malvinas"
$ws.Range("B214").Value = 13
$ws.Range("C214").Value = 0
$ws.Range("D214").Value = 13
$ws.Range("E214").Value = 0
$ws.Range("F214").Value = 0
$ws.Range("G214").Value = 0
$ws.Range("H214").Value = 0
